# Automatische test-sync: 2025-07-23 22:38:50
#
# Appends the 11th test-mail row to the "Logs" sheet, bumps the matching
# "Productinformatie" category count on the "Dashboard" sheet, and extends
# the chart series / conditional-formatting ranges to cover the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append row 21
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A21").Value = "Hebben jullie een handleiding van de EcoPro-500?"
$logs.Range("B21").Value = "mailmind.test@zohomail.eu"
$logs.Range("C21").Value = "Testmail #11: Hebben jullie een handleiding van de EcoPro-500?"
$logs.Range("D21").Value = "Productinformatie"
$logs.Range("E21").Value = "Geachte klant,`nDank u voor uw interesse in de EcoPro-500. Op dit moment hebben wij geen handleiding beschikbaar voor de EcoPro-500. Echter, wij kunnen u wel voorzien van de basisinformatie over het product, zoals specificaties, functies en gebruiksinstructies. Mocht u specifieke vragen hebben of meer gedetailleerde informatie nodig hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F21").Value = "2025-07-23 22:38:27"
$logs.Range("G21").Value = "Ja"
$logs.Range("H21").Value = "Nee"
$logs.Range("I21").Value = "Ja"
$logs.Range("J21").Value = "Nee"

# Writing the multi-line "Antwoord" text auto-expands the row height; put
# it back to the sheet default (matches the other, untouched rows) instead
# of leaving an explicit ht/customHeight override on row 21.
$logs.Rows.Item(21).AutoFit()

# Extend the conditional-formatting blocks so they keep covering the full
# data range (row 2 through the newly-added row 21).
$logs.Range("D2:D20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D21"))
$logs.Range("G2:G20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G21"))
$logs.Range("H2:H20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H21"))
$logs.Range("I2:I20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I21"))
$logs.Range("J2:J20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J21"))

# ---------------------------------------------------------------------
# 2) Dashboard sheet: append the new category tally row 9
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A9").Value = "Productinformatie"
$dash.Range("B9").Value = 1

# ---------------------------------------------------------------------
# 3) Chart on the Dashboard sheet: widen the category/value series refs
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$9,Dashboard!`$B`$2:`$B`$9,1)"
